$d = $word.ActiveDocument

# 1) Change "Student Name" to "Student Names" by inserting a new run "s"
#    right after the existing "Student Name" text (so it ends up as two
#    separate runs, matching the diff).
$found = $d.Content.Find.Execute("Student Name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng = $d.Content.Find.Parent.Duplicate
}

# Use Find to locate "Student Name" and collapse to its end, then insert "s"
$findRange = $d.Content
$findRange.Find.Execute("Student Name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Collapse(0)
$findRange.InsertAfter("s")

# 2) Remove the "_GoBack" bookmark (it no longer appears in the saved doc).
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}
